$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.614.62"

$ws.Range("D3").Value = "1.642.63"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.38"
$ws.Range("E5").Value = "  +1.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  +0.98%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("D12").Value = "1.872.06"
$ws.Range("E12").Value = "  +0.75%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.667.37"
$ws.Range("E13").Value = "  +2.87%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.21"
$ws.Range("E14").Value = "  +3.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.530"
$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.71"
$ws.Range("E16").Value = "  +3.88%  "

$ws.Range("D17").Value = "26.655.58"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").Value = "0.0₃0748"
$ws.Range("E18").Value = "  +0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.20"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.01"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  +2.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.31"
$ws.Range("E22").Value = "  +2.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.53"
$ws.Range("E23").Value = "  +2.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +10.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.03"
$ws.Range("E25").Value = "  -1.25%  "

$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("E28").Value = "  +4.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.79"
$ws.Range("E29").Value = "  +1.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0515"
$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("E31").Value = "  +0.63%  "

$ws.Range("E32").Value = "  +2.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("E33").Value = "  +2.11%  "

$ws.Range("D34").Value = "1.270.79"
$ws.Range("E34").Value = "  +4.78%  "

$ws.Range("E35").Value = "  +2.44%  "

$ws.Range("E36").Value = "  +5.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.41"
$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.530"
$ws.Range("E38").Value = "  +5.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.827"
$ws.Range("E39").Value = "  +2.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  +0.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  +2.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.24"
$ws.Range("E42").Value = "  -1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.47"
$ws.Range("E43").Value = "  +1.22%  "

$ws.Range("D44").Value = "1.782.45"
$ws.Range("E44").Value = "  +0.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.99"
$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.42"
$ws.Range("E46").Value = "  +8.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.59"
$ws.Range("E47").Value = "  +2.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0516"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("E49").Value = "  +2.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("E50").Value = "  +3.52%  "

$ws.Range("E51").Value = "  -0.52%  "
